# Insert a new weekly price record for "Repollo" (Terminal Hortofrutícola
# Agro Chillán) as row 109, pushing the existing rows 109:208 down to
# 110:209 (dimension grows from A1:R208 to A1:R209).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 109..208 down to 110..209 by inserting a blank row at 109.
$ws.Rows("109:109").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A109").Value = 7
$ws.Range("B109").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C109").Value = "Ñuble"
$ws.Range("D109").Value = 44658
$ws.Range("E109").Value = 16
$ws.Range("F109").Value = 100112006
$ws.Range("G109").Value = "Repollo"
$ws.Range("H109").Value = "Crespo record"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 200
$ws.Range("K109").Value = 850
$ws.Range("L109").Value = 900
$ws.Range("M109").Value = 875
$ws.Range("N109").Value = "$/unidad"
$ws.Range("O109").Value = "Provincia de Diguillín"
$ws.Range("P109").Value = 875
$ws.Range("Q109").Value = 1
$ws.Range("R109").Value = "Hortaliza"
